$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4752
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -4142
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4000.7827
$ws.Range("I69").Value = 4128.4287
$ws.Range("K69").Value = 12385.2861
$ws.Range("M69").Value = -11511.2861

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4000.7827
$ws.Range("I72").Value = 4128.4287
$ws.Range("K72").Value = 37155.85830000001
$ws.Range("M72").Value = -32787.85830000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3964.9656
$ws.Range("I74").Value = 3940.4736
$ws.Range("K74").Value = 3940.4736
$ws.Range("M74").Value = -3004.4736

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3044.4707
$ws.Range("I76").Value = 3120
$ws.Range("J76").Value = 2834.6667
$ws.Range("K76").Value = 3120
$ws.Range("L76").Value = 2834.6667
$ws.Range("M76").Value = -2805
$ws.Range("N76").Value = -3464.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3964.9656
$ws.Range("I77").Value = 3940.4736
$ws.Range("K77").Value = 19702.368
$ws.Range("M77").Value = -15022.368

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3044.4707
$ws.Range("I79").Value = 3120
$ws.Range("J79").Value = 2834.6667
$ws.Range("K79").Value = 3120
$ws.Range("L79").Value = 2834.6667
$ws.Range("M79").Value = -2028
$ws.Range("N79").Value = -5018.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2653.3708
$ws.Range("I138").Value = 1540.5883
$ws.Range("J138").Value = 4004.6072
$ws.Range("K138").Value = 4621.7649
$ws.Range("L138").Value = 12013.8216
$ws.Range("M138").Value = 518.2350999999999
$ws.Range("N138").Value = -22293.8216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2831.889
$ws.Range("I132").Value = 1769.4872
$ws.Range("J132").Value = 5594.1333
$ws.Range("K132").Value = 5308.461600000001
$ws.Range("L132").Value = 16782.3999
$ws.Range("M132").Value = -2778.461600000001
$ws.Range("N132").Value = -21842.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 40002796
$ws.Range("I20").Value = 2849.1875
$ws.Range("J20").Value = 111113816
$ws.Range("K20").Value = 2849.1875
$ws.Range("L20").Value = 111113816
$ws.Range("M20").Value = -2602.1875
$ws.Range("N20").Value = -111114310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2156.3333
$ws.Range("I86").Value = 1750
$ws.Range("J86").Value = 2969
$ws.Range("K86").Value = 1750
$ws.Range("L86").Value = 2969
$ws.Range("M86").Value = -627
$ws.Range("N86").Value = -5215

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2156.3333
$ws.Range("I89").Value = 1750
$ws.Range("J89").Value = 2969
$ws.Range("K89").Value = 8750
$ws.Range("L89").Value = 14845
$ws.Range("M89").Value = -3134
$ws.Range("N89").Value = -26077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19999
$ws.Range("J41").Value = 19999
$ws.Range("L41").Value = 19999
$ws.Range("N41").Value = -20855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 19999
$ws.Range("J51").Value = 19999
$ws.Range("L51").Value = 19999
$ws.Range("N51").Value = -21471

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 20360.25
$ws.Range("J59").Value = 20360.25
$ws.Range("L59").Value = 20360.25
$ws.Range("N59").Value = -22650.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 19999
$ws.Range("J61").Value = 19999
$ws.Range("L61").Value = 19999
$ws.Range("N61").Value = -20695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3757.1428
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3757.1428
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 28176.111
$ws.Range("J74").Value = 28176.111
$ws.Range("L74").Value = 28176.111
$ws.Range("N74").Value = -29924.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 28176.111
$ws.Range("J77").Value = 28176.111
$ws.Range("L77").Value = 84528.333
$ws.Range("N77").Value = -93264.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 25363406
$ws.Range("I132").Value = 27778588
$ws.Range("J132").Value = 16668747
$ws.Range("K132").Value = 83335764
$ws.Range("L132").Value = 50006241
$ws.Range("M132").Value = -83333234
$ws.Range("N132").Value = -50011301

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1135.965
$ws.Range("I134").Value = 1022.39215
$ws.Range("J134").Value = 2101.3333
$ws.Range("K134").Value = 3067.17645
$ws.Range("L134").Value = 6303.999899999999
$ws.Range("M134").Value = -532.1764499999999
$ws.Range("N134").Value = -11373.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 939.5
$ws.Range("I68").Value = 1033.3334
$ws.Range("J68").Value = 913.9091
$ws.Range("K68").Value = 3100.0002
$ws.Range("L68").Value = 2741.7273
$ws.Range("M68").Value = -2289.0002
$ws.Range("N68").Value = -4363.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 939.5
$ws.Range("I71").Value = 1033.3334
$ws.Range("J71").Value = 913.9091
$ws.Range("K71").Value = 9300.000599999999
$ws.Range("L71").Value = 8225.1819
$ws.Range("M71").Value = -5244.000599999999
$ws.Range("N71").Value = -16337.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 27778114
$ws.Range("J107").Value = 100000344
$ws.Range("L107").Value = 300001032
$ws.Range("N107").Value = -300004872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 534
$ws.Range("I113").Value = 521.44446
$ws.Range("J113").Value = 546.55554
$ws.Range("K113").Value = 1564.33338
$ws.Range("L113").Value = 1639.66662
$ws.Range("M113").Value = 605.66662
$ws.Range("N113").Value = -5979.66662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4798.129
$ws.Range("I131").Value = 426.42856
$ws.Range("J131").Value = 6073.2085
$ws.Range("K131").Value = 1279.28568
$ws.Range("L131").Value = 18219.6255
$ws.Range("M131").Value = 3760.71432
$ws.Range("N131").Value = -28299.6255

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 56645748
$ws.Range("I80").Value = 101800940
$ws.Range("J80").Value = 201745
$ws.Range("K80").Value = 101800940
$ws.Range("L80").Value = 201745
$ws.Range("M80").Value = -101799942
$ws.Range("N80").Value = -203741

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 56645748
$ws.Range("I83").Value = 101800940
$ws.Range("J83").Value = 201745
$ws.Range("K83").Value = 509004700
$ws.Range("L83").Value = 1008725
$ws.Range("M83").Value = -508999708
$ws.Range("N83").Value = -1018709

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1629.7333
$ws.Range("I102").Value = 1386.6086
$ws.Range("J102").Value = 2428.5715
$ws.Range("K102").Value = 1386.6086
$ws.Range("L102").Value = 2428.5715
$ws.Range("M102").Value = 235.3914
$ws.Range("N102").Value = -5672.5715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2849.7273
$ws.Range("I126").Value = 2617.875
$ws.Range("J126").Value = 3468
$ws.Range("K126").Value = 7853.625
$ws.Range("L126").Value = 10404
$ws.Range("M126").Value = -5383.625
$ws.Range("N126").Value = -15344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 41668948
$ws.Range("I82").Value = 41668948
$ws.Range("K82").Value = 41668948
$ws.Range("M82").Value = -41668587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 41668948
$ws.Range("I85").Value = 41668948
$ws.Range("K85").Value = 41668948
$ws.Range("M85").Value = -41667700

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6727.8945
$ws.Range("I93").Value = 7452.5
$ws.Range("J93").Value = 2863.3333
$ws.Range("K93").Value = 7452.5
$ws.Range("L93").Value = 2863.3333
$ws.Range("M93").Value = -6204.5
$ws.Range("N93").Value = -5359.3333
